$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: "Assim, o utilizador..." paragraph: merge/extend text, add
# firstLine indent, and insert a whole new paragraph ("Posto isto, ...")
# before the page-break run.
# ---------------------------------------------------------------------------
$oldP1P2 = "Assim, o utilizador pode selecionar um dos mapas disponíveis, acedendo depois a um menu de navegação. Com isto, o utilizador pode escolher as diferentes opções de navegação. Em relação a estas opções, o utilizador pode selecionar tanto o ponto de origem como de destino, adicionar ou remover pontos de interesse (pontos que serão visitados no percurso). Além disso, é de esperar que o programa evite zonas inacessíveis, indicando uma rota rápida e correta."

$newP1 = "Assim, o utilizador pode selecionar um dos mapas disponíveis, acedendo depois a um menu de navegação. Com isto, o utilizador pode escolher as diferentes opções de navegação. Em relação a estas opções, o utilizador pode selecionar tanto o ponto de origem como de destino, adicionar ou remover pontos de interesse (pontos que serão visitados no percurso). Além disso, é também dada a opção de remover zonas inacessíveis ou adicioná-las, o que obviamente pode implicar mudanças na conetividade do grafo e alterações drásticas no percurso final."

$newP2 = "Posto isto, em relação ao percurso, o utilizador ainda pode escolher se pretender efetuar um percurso mais econômico, isto é, minimizar o custo da viagem (evitando pórticos) mas como consequência aumentar o tempo de viagem, ou o percurso mais rápido, ou seja, minimizar o tempo de viagem mas aumentar o custo. A opção “Path Visualization” permite então visualizar o percurso de acordo com a informação dada pelo utilizador."

# Replace the whole old text (no paragraph split yet) with the concatenated
# new text of both paragraphs.
$r = $d.Content
$found = $r.Find.Execute($oldP1P2, $true, $false, $false, $false, $false, $true, 1, $false, ($newP1 + $newP2), 2)
Write-Output "Hunk1 replace found=$found"

# Now split into two paragraphs right at the boundary, just after newP1 (and
# before the page-break run), using InsertParagraphAfter so the page-break
# run is left untouched in its own run.
$r2 = $d.Content
$r2.Find.Execute($newP1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Collapse(0)
$r2.InsertParagraphAfter()
Write-Output "Hunk1 paragraph split done"

# Add firstLine indent (708 twips = 35.4pt) to both paragraphs.
$r3 = $d.Content
$r3.Find.Execute($newP1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Paragraphs(1).Range.ParagraphFormat.FirstLineIndent = 35.4

$r4 = $d.Content
$r4.Find.Execute($newP2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r4.Paragraphs(1).Range.ParagraphFormat.FirstLineIndent = 35.4

Write-Output "Hunk1 done"
